$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 1.83
$ws.Range("P2").Value = 1.83
$ws.Range("U2").Value = 1.83
$ws.Range("V2").Value = 1.98
$ws.Range("X2").Value = 1.36

# Row 3 updates
$ws.Range("G3").Value = 1.33
$ws.Range("H3").Value = 4.1
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 1.95
$ws.Range("K3").Value = 2.1
$ws.Range("L3").Value = 12
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 1.17
$ws.Range("U3").Value = 1.53
$ws.Range("V3").Value = 2.38
$ws.Range("W3").Value = 3.4
$ws.Range("Z3").Value = 4.75
$ws.Range("AB3").Value = 7.5
$ws.Range("AC3").Value = 17
$ws.Range("AE3").Value = 6
$ws.Range("AF3").Value = 9.5
$ws.Range("AG3").Value = 41
$ws.Range("AH3").Value = 201
$ws.Range("AI3").Value = 19
$ws.Range("AJ3").Value = 51
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 201
$ws.Range("AM3").Value = 151
$ws.Range("AN3").Value = 151
$ws.Range("AP3").Value = 1.93
$ws.Range("AQ3").Value = 1.93
